$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G9").Value = 0.89700000000000002
$ws.Range("H9").Value = 0.59699999999999998
$ws.Range("I9").Value = 0.878
$ws.Range("J9").Value = 0.61299999999999999
$ws.Range("K9").Value = 0.92400000000000004
$ws.Range("L9").Value = 0.63100000000000001

$ws.Range("T9").Select()
